# Generate Report for Handoff
# Updates the "b.md" row across the Overview, zh-cn and de-de sheets to reflect
# that the file is ready for handoff (a new handoff xliff was generated), and
# records the version-mismatch error detail.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Overview sheet - row 3 is the "b.md" entry
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-11-09 05:54:54"

# ---------------------------------------------------------------------------
# zh-cn sheet - row 3 is the "b.md" entry
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("F3").Value = "'False"
$wsZhCn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-11-09 05:54:40"
$wsZhCn.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/53363d3f7baf83d7ef185ad0f4d65f85ae42f22b/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/67d109a9c14d18788cbd0d9d77424c68095ea304/e2e/b.md."
$wsZhCn.Range("P1").ColumnWidth = 39.14

# ---------------------------------------------------------------------------
# de-de sheet - row 3 is the "b.md" entry
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("F3").Value = "'False"
$wsDeDe.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$wsDeDe.Range("H3").Value = "2016-11-09 05:54:54"
$wsDeDe.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/53363d3f7baf83d7ef185ad0f4d65f85ae42f22b/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/67d109a9c14d18788cbd0d9d77424c68095ea304/e2e/b.md."
$wsDeDe.Range("P1").ColumnWidth = 39.14
